# Commit: "Fri, Jul 24, 2020  7:05:21 AM"
#
# Two logical changes are applied to the deck:
#
#  1. The table on the slide that still used the (legacy) default table
#     style is switched to a different table-style gallery entry
#     (tableStyleId {5F7B5029-E678-424C-8F58-7F5BBA097CB4} ->
#     {66F45CEB-F9C9-40EC-A3E6-34CA9074E7C9}).
#
#  2. The presentation's theme ("Integral") is recoloured to the stock
#     PowerPoint "Office Theme" palette (the 12 theme colours that used
#     to hold the Integral palette now hold the Office palette), which
#     is what the underlying canonical-XML diff shows as theme1.xml's
#     new content.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Re-style the table.
# ---------------------------------------------------------------------
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle("{66F45CEB-F9C9-40EC-A3E6-34CA9074E7C9}")
        }
    }
}

# ---------------------------------------------------------------------
# 2) Recolour the theme from "Integral" to the default "Office Theme"
#    palette.
# ---------------------------------------------------------------------
function RGBVal($r, $g, $b) { return $r + ($g * 256) + ($b * 65536) }

# Index order matches MsoThemeColorSchemeIndex:
#  1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
$officeThemeColors = @(
    (RGBVal 0x00 0x00 0x00),  # dk1      000000
    (RGBVal 0xFF 0xFF 0xFF),  # lt1      FFFFFF
    (RGBVal 0x44 0x54 0x6A),  # dk2      44546A
    (RGBVal 0xE7 0xE6 0xE6),  # lt2      E7E6E6
    (RGBVal 0x5B 0x9B 0xD5),  # accent1  5B9BD5
    (RGBVal 0xED 0x7D 0x31),  # accent2  ED7D31
    (RGBVal 0xA5 0xA5 0xA5),  # accent3  A5A5A5
    (RGBVal 0xFF 0xC0 0x00),  # accent4  FFC000
    (RGBVal 0x44 0x72 0xC4),  # accent5  4472C4
    (RGBVal 0x70 0xAD 0x47),  # accent6  70AD47
    (RGBVal 0x05 0x63 0xC1),  # hlink    0563C1
    (RGBVal 0x95 0x4F 0x72)   # folHlink 954F72
)

$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme
for ($k = 1; $k -le $colorScheme.Count; $k++) {
    $colorScheme.Item($k).RGB = $officeThemeColors[$k - 1]
}
